$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (single dot, e.g. "228.25").
# Cells whose values retain multiple dots (e.g. "37.938.06") are safely
# kept as text automatically and do not need this.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values exactly as they appear in the refreshed
# cryptos listing (coin prices, 1h volume deltas, and the two coin-name /
# link swaps lower in the table).
$ws.Range("D2").Value = '37.938.06'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.035.04'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '228.25'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '60.48'
$ws.Range("E7").Value = '  +2.95%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.379'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").Value = '0.0823'
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '2.335.84'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '14.53'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").Value = '21.40'
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").Value = '0.762'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = '2.043.22'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '37.850.22'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").Value = '69.70'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '5.90'
$ws.Range("E20").Value = '  -7.28%  '
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").Value = '222.53'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '166.80'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '9.30'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '0.130'
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").Value = '18.85'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("E32").Value = '  +8.22%  '
$ws.Range("D33").Value = '4.41'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").Value = '4.51'
$ws.Range("E35").Value = '  -1.97%  '
$ws.Range("D36").Value = '6.37'
$ws.Range("E36").Value = '  +5.05%  '
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.539.33'
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '17.65'
$ws.Range("E41").Value = '  +6.77%  '
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '96.15'
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  -2.93%  '
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("D49").Value = '2.96'
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").Value = '7.08'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").Value = '2.225.45'
$ws.Range("E51").Value = '  -0.93%  '
